$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 597.6180316666668
$ws.Range("H2").Value = 1792.854095
$ws.Range("I2").Value = 0.4787032177461795
$ws.Range("J2").Value = 0.4787032177461795
$ws.Range("M2").Value = 0.119457
$ws.Range("N2").Value = 0.358371
$ws.Range("Q2").Value = 71.389657208805
$ws.Range("R2").Value = 642.5069148792451
$ws.Range("S2").Value = 0.4787032177461795
$ws.Range("T2").Value = 0.4787032177461795

# Row 3
$ws.Range("I3").Value = 0.3987839532217896
$ws.Range("J3").Value = 0.3987839532217896
$ws.Range("M3").Value = 0.119457
$ws.Range("N3").Value = 0.358371
$ws.Range("Q3").Value = 59.47118938308599
$ws.Range("R3").Value = 535.240704447774
$ws.Range("S3").Value = 0.3987839532217896
$ws.Range("T3").Value = 0.3987839532217896

# Row 4
$ws.Range("G4").Value = 152.9462786666667
$ws.Range("H4").Value = 458.838836
$ws.Range("I4").Value = 0.1225128290320309
$ws.Range("J4").Value = 0.1225128290320309
$ws.Range("M4").Value = 0.119457
$ws.Range("N4").Value = 0.358371
$ws.Range("Q4").Value = 18.270503610684
$ws.Range("R4").Value = 164.434532496156
$ws.Range("S4").Value = 0.1225128290320309
$ws.Range("T4").Value = 0.1225128290320309
